$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Step 1: copy/paste formats BEFORE any values change, using the
# worksheet's own existing (pre-edit) cells as format donors. ---

# Row 6 needs to look like the old "College Name / KG RCET" row (style 2/2)
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteFormats) | Out-Null

# Row 7's B cell needs to move from style 2 to style 3 (like B8)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null

# Rows 8-10 already have the correct styles (2/3), nothing to copy.

# Row 11 becomes blank; give it the borderless "gap" look (style 10/12-ish)
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A11").PasteSpecial($xlPasteFormats) | Out-Null

# Row 12 becomes blank; give it the borderless "gap" look too
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B12").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Remove the stray border that tagged along on B11's date format
$ws.Range("A11:B12").Borders.LineStyle = -4142

# Row 13 goes back to a fully default / unformatted state (the row
# effectively disappears from the authored data)
$ws.Range("A13:B13").ClearContents()
$ws.Range("A13:B13").Style = "Normal"

# --- Step 2: now that formatting is in place, update the cell values ---

$ws.Range("A6").Value = "College Name"
$ws.Range("B6").Value = "KG RCET"

$ws.Range("A7").Value = "Degree Type"
$ws.Range("B7").Value = "Engineering"

$ws.Range("A8").Value = "Stream"
$ws.Range("B8").Value = "Electrical and Electronics"

$ws.Range("A9").Value = "Location"
$ws.Range("B9").Value = "Hyderabad"

$ws.Range("A10").Value = "Message"
$ws.Range("B10").Value = "Testing the message box"

$ws.Range("A11").ClearContents()
$ws.Range("B11").ClearContents()

$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()

$ws.Range("A13").ClearContents()
$ws.Range("B13").ClearContents()

# Restore the selection shown in the sheet view
$ws.Range("A1:B10").Select()
$ws.Range("B8").Activate()
